## Regenerate merged AHB files
## - Rename the "_old" / "_new" header columns to the FV2304 / FV2310 vintage
##   names used by the new merge pipeline.
## - Turn the data range into a real Excel Table (Table1) with autofilter.
## - Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (A1:J1 "_old" -> "_FV2304", K1 "diff" stays,
#        L1:U1 "_new" -> "_FV2310") ---------------------------------------
$headers = [ordered]@{
    "A1" = "Segmentname_FV2304"
    "B1" = "Segmentgruppe_FV2304"
    "C1" = "Segment_FV2304"
    "D1" = "Datenelement_FV2304"
    "E1" = "Segment ID_FV2304"
    "F1" = "Code_FV2304"
    "G1" = "Qualifier_FV2304"
    "H1" = "Beschreibung_FV2304"
    "I1" = "Bedingungsausdruck_FV2304"
    "J1" = "Bedingung_FV2304"
    "L1" = "Segmentname_FV2310"
    "M1" = "Segmentgruppe_FV2310"
    "N1" = "Segment_FV2310"
    "O1" = "Datenelement_FV2310"
    "P1" = "Segment ID_FV2310"
    "Q1" = "Code_FV2310"
    "R1" = "Qualifier_FV2310"
    "S1" = "Beschreibung_FV2310"
    "T1" = "Bedingungsausdruck_FV2310"
    "U1" = "Bedingung_FV2310"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- 2) Convert A1:U79 into an Excel Table named "Table1" -----------------
$tableRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row (row 1) --------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
